$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2: Hebrew birthday "תשרי" -> "יא תשרי"
$ws.Range("G2").Value = "יא תשרי"

# AK2: cluster "מרכז" -> "דרום"
$ws.Range("AK2").Value = "דרום"

# AM3: updated phone number, with formatting reset back to the default font
# (border + center/center alignment kept, font reset to workbook default)
$ws.Range("AM3").Value = 549247617
$ws.Range("AM3").Style = "Normal"
$ws.Range("AM3").Borders.Color = 0
$ws.Range("AM3").Borders.LineStyle = 1
$ws.Range("AM3").HorizontalAlignment = -4108
$ws.Range("AM3").VerticalAlignment = -4108

# Scroll/selection state left by the editing session
$ws.Range("AK2").Select()
$ws.Application.ActiveWindow.ScrollColumn = $ws.Range("U1").Column
